# Generate Report for Handback
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet (zh-cn/de-de columns) and on each language sheet's
#   Status cell.
# - Each language sheet now has a "Latest Target File" (hyperlinked, like the
#   Source File Name cell) and a "Latest Handback File" populated, plus a
#   real "Latest Handback DateTime" (replacing the 0001-01-01 placeholder).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$targetFile = "bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4e801975cb7803105c37071f9fea236d961090bb/e2e/bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.md"

$hyperlinkColor = 15570276   # BGR for FF6495ED, matching the existing hyperlink style

# ---------------------------------------------------------------------------
# Overview sheet: both language status cells flip to the handed-back status.
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, $targetFile, $targetFile)
$zh.Range("I2").Font.Underline = $true
$zh.Range("I2").Font.Color = $hyperlinkColor

$zh.Range("J2").Value = "bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.f7393c2a9b97089adabe42e45bd360de021cc78f.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-30 17:09:37"

$zh.Columns("I").ColumnWidth = 39.17
$zh.Columns("J").ColumnWidth = 39.17
$zh.Columns("C").ColumnWidth = 29.16

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $targetUrl, [System.Reflection.Missing]::Value, $targetFile, $targetFile)
$de.Range("I2").Font.Underline = $true
$de.Range("I2").Font.Color = $hyperlinkColor

$de.Range("J2").Value = "bdaeab86-0f22-48e3-bccb-b33bcdc08d2b.f7393c2a9b97089adabe42e45bd360de021cc78f.de-de.xlf"
$de.Range("K2").Value = "2016-08-30 17:09:44"

$de.Columns("I").ColumnWidth = 39.17
$de.Columns("J").ColumnWidth = 39.17
$de.Columns("C").ColumnWidth = 29.16

# ---------------------------------------------------------------------------
# Overview sheet column widths (zh-cn / de-de status columns widen to fit
# the longer status text).
# ---------------------------------------------------------------------------
$ov.Columns("E").ColumnWidth = 29.16
$ov.Columns("F").ColumnWidth = 29.16
